$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B ("Betrag" -> shifts to column C)
$ws.Columns.Item(2).Insert()

# The inserted column inherited formatting from the old column B (s="1" from A-column
# bleed-through in some cases) - make sure it is unformatted/General like the source diff.
$ws.Columns.Item(2).ClearFormats()

# 2. Populate the new "Belegnr." column (header + values)
$ws.Range("B1").Value = "Belegnr."
$ws.Range("B2").Value = 6000001
$ws.Range("B3").Value = 6000002
$ws.Range("B4").Value = 6000003
$ws.Range("B5").Value = 6000004
$ws.Range("B6").Value = 6000005
$ws.Range("B7").Value = 6000006
$ws.Range("B8").Value = 6000007

# 3. Append a new row (row 9) with date / Belegnr / Betrag, matching formatting of row 8
$ws.Range("A8:C8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("A9").Value = 43160
$ws.Range("B9").Value = 6000011
$ws.Range("C9").Value = 250

$excel.CutCopyMode = 0

# 4. Update the active selection to match the saved workbook state
$ws.Range("C10").Select()

# 5. Register the DATA19 defined name pointing at Blatt1!$I$6:$I$276
$wb.Names.Add('DATA19', '=Blatt1!$I$6:$I$276')
